$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C10 value from 18 to 1 (numeric)
$ws.Range("C10").Value = 1

$wb.Save()
